# Worked on temporal resolution
# Extend the "Demand" sheet (t / EU27.Elec) from a single timestep to a full
# 12-step series, fix up the year-1 demand value, widen the value column,
# and leave "Demand" as the active/selected sheet (matching the commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# New demand values for t = 1..12 (t = 0 stays at 0)
$demandValue = 376098958

for ($t = 1; $t -le 12; $t++) {
    $row = $t + 2
    $ws.Cells.Item($row, 1).Value = $t
    $ws.Cells.Item($row, 2).Value = $demandValue
}

# Column B (EU27.Elec) is now sized to fit the new values
$ws.Columns.Item(2).ColumnWidth = 9.17

# Make "Demand" the active sheet/tab and set the new selection
$ws.Activate() | Out-Null
$ws.Range("D12").Select() | Out-Null
